# Append 7 new daily COVID overview rows (2022-01-31 .. 2022-02-06) to Sheet1,
# directly below the existing last row (537), matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRows = @(
    @("2022-01-31", "overview", "K02000001", "United Kingdom", 17315893, 92368, 51, 155754),
    @("2022-02-01", "overview", "K02000001", "United Kingdom", 17428345, 112458, 219, 156875),
    @("2022-02-02", "overview", "K02000001", "United Kingdom", 17515199, 88085, 534, 157409),
    @("2022-02-03", "overview", "K02000001", "United Kingdom", 17607832, 88171, 303, 157730),
    @("2022-02-04", "overview", "K02000001", "United Kingdom", 17689885, 84053, 254, 157984),
    @("2022-02-05", "overview", "K02000001", "United Kingdom", 17749999, 60578, 259, 158243),
    @("2022-02-06", "overview", "K02000001", "United Kingdom", 17803325, 54095, 75, 158318)
)

$startRow = 538

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Column A holds the date as plain text (matches existing rows, which are
    # inline strings, not Excel date serials), so force text with a leading
    # apostrophe-free NumberFormat of @ and assign as string.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
